$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.672.91'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '3.346.42'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  -0.03%  '
$c = $ws.Range("D5")
$c.Value = "'259.10"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$c = $ws.Range("D6")
$c.Value = "'648.01"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.89%  '
$ws.Range("E7").Value = '  +9.94%  '
$c = $ws.Range("D8")
$c.Value = "'0.460"
$c.Style = "Normal"
$ws.Range("E8").Value = '  +16.60%  '
$ws.Range("E9").Value = '  +23.54%  '
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '3.343.93'
$ws.Range("E11").Value = '  -0.65%  '
$c = $ws.Range("D12")
$c.Value = "'0.209"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +4.75%  '
$c = $ws.Range("D13")
$c.Value = "'43.71"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +21.23%  '
$c = $ws.Range("D14")
$c.Value = "'0.0000267"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +7.41%  '
$ws.Range("D15").Value = '99.202.98'
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("D16").Value = '3.986.36'
$ws.Range("E16").Value = '  +0.63%  '
$c = $ws.Range("D17")
$c.Value = "'5.57"
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").Value = '3.374.28'
$ws.Range("E18").Value = '  +0.23%  '
$c = $ws.Range("D19")
$c.Value = "'7.55"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +21.72%  '
$c = $ws.Range("D20")
$c.Value = "'16.89"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +10.75%  '
$c = $ws.Range("D21")
$c.Value = "'537.03"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +8.87%  '
$c = $ws.Range("D22")
$c.Value = "'3.56"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("E23").Value = '  +8.81%  '
$c = $ws.Range("D24")
$c.Value = "'0.0000212"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.43%  '
$c = $ws.Range("D25")
$c.Value = "'0.438"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +53.09%  '
$c = $ws.Range("D26")
$c.Value = "'103.31"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +16.34%  '
$ws.Range("E27").Value = '  +9.53%  '
$ws.Range("E28").Value = '  +5.66%  '
$ws.Range("D29").Value = '3.527.82'
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("E30").Value = '  +8.08%  '
$ws.Range("E31").Value = '  +0.22%  '
$c = $ws.Range("D32")
$c.Value = "'10.96"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +14.44%  '
$c = $ws.Range("D33")
$c.Value = "'0.191"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.54%  '
$ws.Range("E34").Value = '  +0.42%  '
$c = $ws.Range("D35")
$c.Value = "'29.30"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +5.03%  '
$c = $ws.Range("D36")
$c.Value = "'0.541"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +16.52%  '
$c = $ws.Range("D37")
$c.Value = "'7.75"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +5.68%  '
$c = $ws.Range("D38")
$c.Value = "'2.09"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +5.80%  '
$ws.Range("E39").Value = '  +2.68%  '
$c = $ws.Range("D40")
$c.Value = "'518.20"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.82%  '
$c = $ws.Range("D41")
$c.Value = "'24.71"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.57%  '
$c = $ws.Range("D42")
$c.Value = "'1.31"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +3.13%  '
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D44")
$c.Value = "'0.0432"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +31.57%  '
$c = $ws.Range("D45")
$c.Value = "'0.823"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +5.21%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D46")
$c.Value = "'3.36"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D48")
$c.Value = "'7.88"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +19.65%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D49")
$c.Value = "'5.07"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +8.78%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D50")
$c.Value = "'2.03"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +3.67%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D51")
$c.Value = "'164.15"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.32%  '
